$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.848.60'
$ws.Range('E2').Value = '  -0.91%  '
$ws.Range('D3').Value = '1.561.41'
$ws.Range('E3').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '205.69'
$ws.Range('E5').Value = '  -0.37%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.485'
$ws.Range('E6').Value = '  -0.38%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '21.73'
$ws.Range('E8').Value = '  -1.58%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.246'
$ws.Range('E9').Value = '  -0.19%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0585'
$ws.Range('E10').Value = '  -0.79%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0864'
$ws.Range('E11').Value = '  +0.00%  '
$ws.Range('D12').Value = '1.785.34'
$ws.Range('E12').Value = '  +0.10%  '
$ws.Range('D13').Value = '1.563.68'
$ws.Range('E13').Value = '  +0.62%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '3.72'
$ws.Range('E14').Value = '  -0.83%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.515'
$ws.Range('E15').Value = '  -0.37%  '
$ws.Range('D16').Value = '26.868.30'
$ws.Range('E16').Value = '  -0.84%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '61.28'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '215.23'
$ws.Range('E18').Value = '  +0.34%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.37'
$ws.Range('E19').Value = '  +1.94%  '
$ws.Range('D20').Value = '0.0₃0682'
$ws.Range('E20').Value = '  -0.21%  '
$ws.Range('E21').Value = '  -0.04%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.13'
$ws.Range('E22').Value = '  +0.46%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.19'
$ws.Range('E23').Value = '  -1.24%  '
$ws.Range('E24').Value = '  +1.44%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '153.91'
$ws.Range('E25').Value = '  +1.50%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '6.68'
$ws.Range('E26').Value = '  +1.06%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '14.94'
$ws.Range('E27').Value = '  +0.23%  '
$ws.Range('E29').Value = '  -0.75%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0464'
$ws.Range('E30').Value = '  +0.87%  '
$ws.Range('E31').Value = '  -3.27%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.16'
$ws.Range('E32').Value = '  +0.17%  '
$ws.Range('D33').Value = '1.390.79'
$ws.Range('E33').Value = '  +0.43%  '
$ws.Range('E34').Value = '  -0.18%  '
$ws.Range('E35').Value = '  -1.31%  '
$ws.Range('E36').Value = '  -0.48%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.917'
$ws.Range('E37').Value = '  -2.59%  '
$ws.Range('E38').Value = '  -0.49%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.529'
$ws.Range('E39').Value = '  +2.72%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.811'
$ws.Range('E40').Value = '  +0.19%  '
$ws.Range('E41').Value = '  -0.05%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.991'
$ws.Range('E42').Value = '  +0.12%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.52'
$ws.Range('E43').Value = '  +4.78%  '
$ws.Range('E44').Value = '  +0.90%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.19'
$ws.Range('E45').Value = '  +1.33%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '63.60'
$ws.Range('E46').Value = '  +0.49%  '
$ws.Range('D47').Value = '1.698.48'
$ws.Range('E47').Value = '  +0.21%  '
$ws.Range('E48').Value = '  +1.58%  '
$ws.Range('E49').Value = '  +2.71%  '
$ws.Range('D50').Value = '0.0₇0984'
$ws.Range('E50').Value = '  -0.02%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0952'
$ws.Range('E51').Value = '  +1.09%  '
